# Deploying to gh-pages — add the 2023 data column (T) to the 1.a.2 table
# and tidy up row heights / selection state to match the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column T with the 2023 figures, copying the formatting from
#     the corresponding cell in column S (the previous last year column).
$pairs = @(
    @{ Cell = "T4"; Source = "S4"; Value = 2023 },
    @{ Cell = "T5"; Source = "S5"; Value = 43.1 },
    @{ Cell = "T6"; Source = "S6"; Value = 19.7 },
    @{ Cell = "T7"; Source = "S7"; Value = 7.8 },
    @{ Cell = "T8"; Source = "S8"; Value = 15.6 }
)

foreach ($pair in $pairs) {
    $ws.Range($pair.Source).Copy()
    $ws.Range($pair.Cell).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($pair.Cell).Value = $pair.Value
}

$excel.CutCopyMode = $false

# --- Row height tweaks for the refreshed layout
$ws.Rows.Item(1).RowHeight = 57
$ws.Rows.Item(4).RowHeight = 16.5

# --- Reset the saved selection back to the top-left cell
$ws.Range("A1").Select() | Out-Null
